# Auto-generated edit script applying market-data value updates
# across the Coeurl_Profits workbook (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 833356.9399999999
$ws.Range("I6").Value = 1000019.9
$ws.Range("K6").Value = 3000059.7
$ws.Range("M6").Value = -2999947.7
$ws.Range("H132").Value = 3472.2222
$ws.Range("I132").Value = 2387.5
$ws.Range("K132").Value = 7162.5
$ws.Range("M132").Value = -4632.5
$ws.Range("H133").Value = 77112.5
$ws.Range("J133").Value = 77112.5
$ws.Range("L133").Value = 77112.5
$ws.Range("N133").Value = -87232.5
$ws.Range("H137").Value = 1523.25
$ws.Range("I137").Value = 1335.25
$ws.Range("J137").Value = 1899.25
$ws.Range("K137").Value = 4005.75
$ws.Range("L137").Value = 5697.75
$ws.Range("M137").Value = -1455.75
$ws.Range("N137").Value = -10797.75
$ws.Range("H138").Value = 3409.7974
$ws.Range("I138").Value = 1783.4073
$ws.Range("K138").Value = 5350.2219
$ws.Range("M138").Value = -210.2219000000005

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3520.75
$ws.Range("I122").Value = 3399.4583
$ws.Range("K122").Value = 10198.3749
$ws.Range("M122").Value = -7748.374899999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 406527.3
$ws.Range("I22").Value = 678.5625
$ws.Range("J22").Value = 541810.25
$ws.Range("K22").Value = 678.5625
$ws.Range("L22").Value = 541810.25
$ws.Range("M22").Value = -505.5625
$ws.Range("N22").Value = -542156.25
$ws.Range("H132").Value = 89791
$ws.Range("J132").Value = 89791
$ws.Range("L132").Value = 89791
$ws.Range("N132").Value = -99911
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20112.428
$ws.Range("I31").Value = 29435.723
$ws.Range("K31").Value = 29435.723
$ws.Range("M31").Value = -29140.723
$ws.Range("H34").Value = 20112.428
$ws.Range("I34").Value = 29435.723
$ws.Range("K34").Value = 29435.723
$ws.Range("M34").Value = -29233.723
$ws.Range("H105").Value = 2147.6667
$ws.Range("I105").Value = 1721.5
$ws.Range("K105").Value = 1721.5
$ws.Range("M105").Value = 25.5
$ws.Range("H132").Value = 5328.8164
$ws.Range("I132").Value = 3655.1667
$ws.Range("K132").Value = 10965.5001
$ws.Range("M132").Value = -8435.500100000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 30399096
$ws.Range("I4").Value = 38422052
$ws.Range("J4").Value = 313011.75
$ws.Range("K4").Value = 115266156
$ws.Range("L4").Value = 939035.25
$ws.Range("M4").Value = -115266044
$ws.Range("N4").Value = -939259.25
$ws.Range("H22").Value = 4349.3076
$ws.Range("J22").Value = 7413.5713
$ws.Range("L22").Value = 22240.7139
$ws.Range("N22").Value = -22578.7139
$ws.Range("H27").Value = 4349.3076
$ws.Range("J27").Value = 7413.5713
$ws.Range("L27").Value = 22240.7139
$ws.Range("N27").Value = -22444.7139
$ws.Range("H82").Value = 27784.928
$ws.Range("I82").Value = 18999
$ws.Range("K82").Value = 56997
$ws.Range("M82").Value = -56591
$ws.Range("H85").Value = 27784.928
$ws.Range("I85").Value = 18999
$ws.Range("K85").Value = 56997
$ws.Range("M85").Value = -55593
$ws.Range("H129").Value = 1025.7142
$ws.Range("J129").Value = 1850
$ws.Range("L129").Value = 5550
$ws.Range("N129").Value = -15550
$ws.Range("H131").Value = 23340.299
$ws.Range("J131").Value = 2448.5945
$ws.Range("L131").Value = 7345.7835
$ws.Range("N131").Value = -17425.7835
$ws.Range("H137").Value = 3362.3845
$ws.Range("J137").Value = 5906.6
$ws.Range("L137").Value = 17719.8
$ws.Range("N137").Value = -27919.8

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 22244212
$ws.Range("I18").Value = 37038690
$ws.Range("J18").Value = 52499.5
$ws.Range("K18").Value = 37038690
$ws.Range("L18").Value = 52499.5
$ws.Range("M18").Value = -37038397
$ws.Range("N18").Value = -53085.5
$ws.Range("H22").Value = 3198.1428
$ws.Range("I22").Value = 1537.8572
$ws.Range("J22").Value = 4858.4287
$ws.Range("K22").Value = 1537.8572
$ws.Range("L22").Value = 4858.4287
$ws.Range("M22").Value = -1008.8572
$ws.Range("N22").Value = -5916.4287
$ws.Range("H97").Value = 1312.4
$ws.Range("I97").Value = 2077.3333
$ws.Range("J97").Value = 165
$ws.Range("K97").Value = 2077.3333
$ws.Range("L97").Value = 165
$ws.Range("M97").Value = -1581.3333
$ws.Range("N97").Value = -1157

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1039.9286
$ws.Range("I16").Value = 1152.7391
$ws.Range("J16").Value = 521
$ws.Range("K16").Value = 1152.7391
$ws.Range("L16").Value = 521
$ws.Range("M16").Value = -982.7391
$ws.Range("N16").Value = -861
$ws.Range("H22").Value = 2204.3914
$ws.Range("I22").Value = 1916.8334
$ws.Range("J22").Value = 2305.8823
$ws.Range("K22").Value = 1916.8334
$ws.Range("L22").Value = 2305.8823
$ws.Range("M22").Value = -1621.8334
$ws.Range("N22").Value = -2895.8823
$ws.Range("H27").Value = 2204.3914
$ws.Range("I27").Value = 1916.8334
$ws.Range("J27").Value = 2305.8823
$ws.Range("K27").Value = 1916.8334
$ws.Range("L27").Value = 2305.8823
$ws.Range("M27").Value = -1809.8334
$ws.Range("N27").Value = -2519.8823
$ws.Range("H29").Value = 3909.4
$ws.Range("J29").Value = 3886.75
$ws.Range("L29").Value = 3886.75
$ws.Range("N29").Value = -4476.75
$ws.Range("H132").Value = 6544.4546
$ws.Range("J132").Value = 6998.1665
$ws.Range("L132").Value = 20994.4995
$ws.Range("N132").Value = -26054.4995
$ws.Range("H136").Value = 5443.55
$ws.Range("I136").Value = 4804.7334
$ws.Range("J136").Value = 7360
$ws.Range("K136").Value = 14414.2002
$ws.Range("L136").Value = 22080
$ws.Range("M136").Value = -11864.2002
$ws.Range("N136").Value = -27180

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 460
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H96").Value = 2526.1428
$ws.Range("I96").Value = 2747.6667
$ws.Range("J96").Value = 2360
$ws.Range("K96").Value = 2747.6667
$ws.Range("L96").Value = 2360
$ws.Range("M96").Value = -1374.6667
$ws.Range("N96").Value = -5106
$ws.Range("H137").Value = 100694.5
$ws.Range("J137").Value = 100694.5
$ws.Range("L137").Value = 100694.5
$ws.Range("N137").Value = -110894.5
